$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Starting layout (before):
#   Row1: headers
#   Row2: TC001 - unauthenticated user - "some keys" search (style row A)
#   Row3: TC002 - authenticated user   - "some keys" search (style row B)
#
# Target layout (after):
#   Row1: headers (unchanged)
#   Row2: TC001 step 1 - unauthenticated user - entering VALID keys
#   Row3: TC001 step 2 - unauthenticated user - entering INVALID keys
#   Row4: TC002 step 1 - authenticated user   - entering VALID keys
#   Row5: TC002 step 2 - authenticated user   - entering INVALID keys
# ---------------------------------------------------------------------------

# 1) Make room: insert a fresh row under row 2 (for the new "invalid keys"
#    unauthenticated case) by copying row 2's formatting down into it.
$ws.Rows("3:3").Insert()
$ws.Range("A2:H2").Copy($ws.Range("A3:H3"))

# 2) Make room: insert a fresh row under what is now row 4 (old TC002 row,
#    pushed down by the previous insert) for the new "invalid keys"
#    authenticated case, copying that row's formatting down into it.
$ws.Rows("5:5").Insert()
$ws.Range("A4:H4").Copy($ws.Range("A5:H5"))

# ---------------------------------------------------------------------------
# Row 2: TC001 - unauthenticated user entering VALID keys
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "As an unauthenticated user,  entering  valid keys to search related products."
$ws.Range("C2").Value = "Products related to the valid search keys will displayed properly and don't show unrelated products."

# ---------------------------------------------------------------------------
# Row 3: TC001 - unauthenticated user entering INVALID keys
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Entering  invalid keys to search related products."
$ws.Range("C3").Value = "The system shows the message: ""Không tìm thấy kết quả phù hợp!"""

# ---------------------------------------------------------------------------
# Row 4: TC002 - authenticated user entering VALID keys
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "As an authenticated user,  entering valid keys to search related products."
$ws.Range("C4").Value = "Products related to the valid search keys will displayed properly and don't show unrelated products."

# ---------------------------------------------------------------------------
# Row 5: TC002 - authenticated user entering INVALID keys
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Entering invalid keys to search related products."
$ws.Range("C5").Value = "The system shows the message: ""Không tìm thấy kết quả phù hợp!"""

# ---------------------------------------------------------------------------
# All four data rows now use a shorter row height.
# ---------------------------------------------------------------------------
$ws.Range("A2:H5").RowHeight = 30

# ---------------------------------------------------------------------------
# Final cursor position, as left by the author.
# ---------------------------------------------------------------------------
[void]$ws.Range("F6").Select()
